function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws "D2" "62.199.57"
Set-TextCell $ws "E2" "  +0.74%  "

Set-TextCell $ws "D3" "2.415.03"
Set-TextCell $ws "E3" "  +1.21%  "

Set-TextCell $ws "E4" "  +0.08%  "

Set-TextCell $ws "D5" "562.52"
Set-TextCell $ws "E5" "  +1.84%  "

Set-TextCell $ws "D6" "142.87"
Set-TextCell $ws "E6" "  +0.47%  "

Set-TextCell $ws "E7" "  +0.23%  "

Set-TextCell $ws "D8" "0.530"
Set-TextCell $ws "E8" "  +1.40%  "

Set-TextCell $ws "D9" "2.411.50"
Set-TextCell $ws "E9" "  +1.11%  "

Set-TextCell $ws "E10" "  +0.75%  "

Set-TextCell $ws "E11" "  -2.16%  "

Set-TextCell $ws "D12" "5.32"
Set-TextCell $ws "E12" "  -0.89%  "

Set-TextCell $ws "D13" "0.351"
Set-TextCell $ws "E13" "  -0.51%  "

Set-TextCell $ws "D14" "25.66"
Set-TextCell $ws "E14" "  -0.96%  "

Set-TextCell $ws "E15" "  +1.01%  "

Set-TextCell $ws "D16" "2.851.86"
Set-TextCell $ws "E16" "  +1.43%  "

Set-TextCell $ws "D17" "62.107.63"
Set-TextCell $ws "E17" "  +1.10%  "

Set-TextCell $ws "D18" "2.411.39"
Set-TextCell $ws "E18" "  +1.19%  "

Set-TextCell $ws "D19" "11.29"
Set-TextCell $ws "E19" "  +1.73%  "

Set-TextCell $ws "B20" "BitcoinCash"
Set-TextCell $ws "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws "D20" "323.36"
Set-TextCell $ws "E20" "  +0.18%  "

Set-TextCell $ws "B21" "Polkadot"
Set-TextCell $ws "C21" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D21" "4.17"
Set-TextCell $ws "E21" "  -0.03%  "

Set-TextCell $ws "D22" "6.80"
Set-TextCell $ws "E22" "  +1.69%  "

Set-TextCell $ws "E23" "  -0.12%  "

Set-TextCell $ws "D24" "65.65"
Set-TextCell $ws "E24" "  +1.64%  "

Set-TextCell $ws "D25" "1.71"
Set-TextCell $ws "E25" "  -2.50%  "

Set-TextCell $ws "D26" "9.02"
Set-TextCell $ws "E26" "  -1.13%  "

Set-TextCell $ws "D27" "577.98"
Set-TextCell $ws "E27" "  +5.65%  "

Set-TextCell $ws "D28" "0.0₃0949"
Set-TextCell $ws "E28" "  +3.16%  "

Set-TextCell $ws "D29" "2.530.80"
Set-TextCell $ws "E29" "  +2.31%  "

Set-TextCell $ws "E30" "  +0.41%  "

Set-TextCell $ws "D31" "8.22"
Set-TextCell $ws "E31" "  -1.05%  "

Set-TextCell $ws "D32" "1.43"
Set-TextCell $ws "E32" "  +0.78%  "

Set-TextCell $ws "E33" "  +0.33%  "

Set-TextCell $ws "D34" "1.87"
Set-TextCell $ws "E34" "  +1.01%  "

Set-TextCell $ws "D35" "1.53"
Set-TextCell $ws "E35" "  -0.23%  "

Set-TextCell $ws "E36" "  +0.24%  "

Set-TextCell $ws "D37" "5.55"
Set-TextCell $ws "E37" "  -3.75%  "

Set-TextCell $ws "D38" "4.70"
Set-TextCell $ws "E38" "  -1.04%  "

Set-TextCell $ws "B39" "PolygonEcosystemToken"
Set-TextCell $ws "C39" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell $ws "D39" "0.381"
Set-TextCell $ws "E39" "  -0.13%  "

Set-TextCell $ws "B40" "Monero"
Set-TextCell $ws "C40" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D40" "152.32"
Set-TextCell $ws "E40" "  +4.38%  "

Set-TextCell $ws "D41" "18.67"
Set-TextCell $ws "E41" "  +0.49%  "

Set-TextCell $ws "D42" "1.80"
Set-TextCell $ws "E42" "  -6.76%  "

Set-TextCell $ws "D43" "0.995"
Set-TextCell $ws "E43" "  -0.51%  "

Set-TextCell $ws "D44" "2.29"
Set-TextCell $ws "E44" "  +0.51%  "

Set-TextCell $ws "D45" "148.37"
Set-TextCell $ws "E45" "  -0.27%  "

Set-TextCell $ws "D46" "3.64"
Set-TextCell $ws "E46" "  +0.53%  "

Set-TextCell $ws "D47" "0.0534"
Set-TextCell $ws "E47" "  +0.72%  "

Set-TextCell $ws "D48" "20.07"
Set-TextCell $ws "E48" "  -1.32%  "

Set-TextCell $ws "D49" "0.593"
Set-TextCell $ws "E49" "  +1.47%  "

Set-TextCell $ws "D50" "0.0916"
Set-TextCell $ws "E50" "  +1.00%  "

Set-TextCell $ws "D51" "0.0227"
Set-TextCell $ws "E51" "  +1.30%  "
